# Swap the contents of columns B:AC between each of the following row pairs.
# Column A (the running index) stays put on its own row; only the data
# (id .. PL_AhUnder) swaps places between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data spans columns B (2) through AC (29)
$firstCol = 2
$lastCol = 29

$rowPairs = @(
    @(15, 16),
    @(18, 19),
    @(27, 28),
    @(29, 30),
    @(42, 43),
    @(68, 69),
    @(70, 71),
    @(78, 79),
    @(81, 82),
    @(87, 88),
    @(90, 91),
    @(92, 93),
    @(96, 97),
    @(100, 101),
    @(112, 113),
    @(114, 115)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
